$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Program values for row 2 and row 3 (physics courses)
$ws.Range("A2").Value = "TUM_LMU_QUANTUM_SCIENCE_TECHNOLOGY"
$ws.Range("A3").Value = "TUM_PHYSICS_NUCLEAR"

# Remove the rest of the former program rows (4-7) entirely
$ws.Range("A4:B7").Clear()

# Shrink the data validation list range from B1:B7 down to B1:B3
# by removing validation from the now-empty rows
$ws.Range("B4:B7").Validation.Delete()

# Trim the trailing used-range rows (996-1000) that are no longer present
$ws.Range("A996:B1000").EntireRow.Delete()
